# Clean up the cages worksheet data: remove invalid/extra rows and
# keep only the valid, fully-populated rows 2-6, fixing the Material
# values (wood -> Wood) along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything below the header row first, then the whole used range
# will be rebuilt with only the valid rows.
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "CageNumber"
$ws.Range("B1").Value = "length"
$ws.Range("C1").Value = "height"
$ws.Range("D1").Value = "width"
$ws.Range("E1").Value = "Material"

# Data rows (only the valid ones survive the cleanup)
$data = @(
    @(1, 20, 30, 40, "Wood"),
    @(2, 30, 40, 25, "Plastic"),
    @(3, 60, 50, 60, "Metal"),
    @(4, 20, 30, 40, "Wood"),
    @(5, 15, 30, 25, "Plastic")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Column widths to match the final layout
$ws.Range("A:A").ColumnWidth = 13.42578125
$ws.Range("B:C").ColumnWidth = 12.140625
$ws.Range("D:D").ColumnWidth = 12.5703125
$ws.Range("E:E").ColumnWidth = 13.85546875

# Selection / view tweaks
$ws.Range("A1:E1").Select()

$wb.Save()
